# Mosaic cheatsheet edit:
#   "    mean(width) "  ->  "    mean(width)) "
# i.e. the `summarize(mean_width = mean(width) )` line gets an extra
# closing paren (e.g. wrapping it in another call), turning the trailing
# run "mean(width) " into two runs "mean(" + "width)) ".

$p = $ppt.ActivePresentation

# Locate the shape that contains the target code snippet, searching every
# slide instead of hard-coding indices so the script is resilient to the
# exact slide/shape numbering.
$targetShape = $null
$searchText = "mean(width)"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.Contains($searchText)) {
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -eq $null) {
    Write-Output "Could not locate shape containing 'mean(width)'; no changes made."
} else {
    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text

    $pos0 = $full.IndexOf($searchText)      # 0-based offset of "mean(width)"
    $start1 = $pos0 + 1                     # 1-based offset (TextRange is 1-based)

    # Split point matches the original run boundary used by the authored
    # edit: "mean(" (5 chars) stays untouched, "width) " (7 chars,
    # including the trailing space) is replaced by "width)) " so the
    # visible text becomes "mean(width)) ".
    $tailStart = $start1 + 5
    $tailLen = 7
    $tailRange = $tr.Characters($tailStart, $tailLen)

    if ($tailRange.Text -eq "width) ") {
        $tailRange.Text = "width)) "
    } else {
        # Fallback: if the surrounding text does not match exactly what we
        # expect, fall back to a straightforward whole-match replace so the
        # edit still lands.
        $wholeRange = $tr.Characters($start1, $searchText.Length)
        $wholeRange.Text = "mean(width))"
    }
}
